$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7553
$ws.Range("K3").Value = 7803
$ws.Range("C4").Value = 1833
$ws.Range("G4").Value = 1473
$ws.Range("K4").Value = 1643
$ws.Range("K5").Value = 556
$ws.Range("K6").Value = 8701
$ws.Range("C7").Value = 27846
$ws.Range("G7").Value = 24383
$ws.Range("K7").Value = 26256

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 228
$ws.Range("K4").Value = 93
$ws.Range("K6").Value = 192
$ws.Range("K7").Value = 790
$ws.Range("K8").Value = 1716
$ws.Range("K10").Value = 159
$ws.Range("K14").Value = 124
$ws.Range("K16").Value = 62
$ws.Range("K17").Value = 49
$ws.Range("K20").Value = 641
$ws.Range("K23").Value = 263
$ws.Range("K25").Value = 121
$ws.Range("K27").Value = 253
$ws.Range("K29").Value = 1443
$ws.Range("K31").Value = 309
$ws.Range("K33").Value = 1104
$ws.Range("K34").Value = 151
$ws.Range("K37").Value = 867
$ws.Range("K43").Value = 217
$ws.Range("K44").Value = 211
$ws.Range("K45").Value = 38
$ws.Range("K51").Value = 340
$ws.Range("K52").Value = 675
$ws.Range("K53").Value = 329
$ws.Range("K57").Value = 106
$ws.Range("C63").Value = 277
$ws.Range("G63").Value = 284
$ws.Range("K63").Value = 77
$ws.Range("K65").Value = 614
$ws.Range("K67").Value = 1024
$ws.Range("K73").Value = 235
$ws.Range("K75").Value = 84
$ws.Range("K76").Value = 362
$ws.Range("K79").Value = 644
$ws.Range("K80").Value = 99
$ws.Range("K83").Value = 558
$ws.Range("K84").Value = 214
$ws.Range("K85").Value = 1203
$ws.Range("K88").Value = 280
$ws.Range("K90").Value = 253
$ws.Range("K91").Value = 311
$ws.Range("K94").Value = 351
$ws.Range("K96").Value = 279
$ws.Range("K97").Value = 216
$ws.Range("K99").Value = 440
$ws.Range("C101").Value = 27846
$ws.Range("G101").Value = 24383
$ws.Range("K101").Value = 26256

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 85
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 249
$ws.Range("K7").Value = 790

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 397
$ws.Range("K3").Value = 416
$ws.Range("K7").Value = 1203

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 187
$ws.Range("K7").Value = 675

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 86
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 329

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 480
$ws.Range("K3").Value = 518
$ws.Range("K6").Value = 574
$ws.Range("K7").Value = 1716

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 194
$ws.Range("K5").Value = 13
$ws.Range("K6").Value = 132
$ws.Range("K7").Value = 558

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 276
$ws.Range("K3").Value = 392
$ws.Range("K6").Value = 353
$ws.Range("K7").Value = 1104

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K6").Value = 263
$ws.Range("K7").Value = 867

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 228
$ws.Range("K7").Value = 614

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 118
$ws.Range("K7").Value = 440

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 93
$ws.Range("K7").Value = 309

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K5").Value = 25
$ws.Range("K6").Value = 293
$ws.Range("K7").Value = 1024

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K3").Value = 87
$ws.Range("K7").Value = 214

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K4").Value = 67
$ws.Range("K7").Value = 1443

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 211

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 80
$ws.Range("K7").Value = 362

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 75
$ws.Range("K7").Value = 192

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K2").Value = 48
$ws.Range("K7").Value = 159

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 91
$ws.Range("K7").Value = 263

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 145
$ws.Range("K5").Value = 8
$ws.Range("K7").Value = 311

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K6").Value = 168
$ws.Range("K7").Value = 644

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 200
$ws.Range("K6").Value = 184
$ws.Range("K7").Value = 641

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 83
$ws.Range("K7").Value = 351

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K2").Value = 49
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 235

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 69
$ws.Range("K7").Value = 228

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 120
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 74
$ws.Range("K3").Value = 88
$ws.Range("K7").Value = 280

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 67
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 253

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 93
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 253

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 93
$ws.Range("K3").Value = 94
$ws.Range("K7").Value = 340

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 62
